# change to tropo tau
# Adds a "# UT Samples" column to both the "Not Sorted" (A:D/E) block and
# the "Sorted By Increasing Mode" (G:J/K) block on both worksheets.

$wb = $excel.ActiveWorkbook

# Flight -> number of UT samples (same lookup table used on both sheets)
$samples = @{
    "RF05" = 96
    "RF06" = 39
    "RF07" = 44
    "RF08" = 48
    "RF09" = 48
    "RF10" = 75
    "RF11" = 54
    "RF12" = 128
    "RF13" = 86
    "RF14" = 47
}

foreach ($ws in @($wb.Worksheets.Item(1), $wb.Worksheets.Item(2))) {

    # Insert a new blank column at E; this shifts the old E:J block to F:K
    # and carries over formatting (font/border/alignment) from row 2-13.
    $ws.Range("E:E").Insert()

    # --- Header row (row 3) ---
    $ws.Range("E3").Value = "# UT Samples"
    $ws.Range("K3").Value = "# UT Samples"
    $ws.Range("E3").NumberFormat = "General"
    $ws.Range("K3").NumberFormat = "General"
    $ws.Range("K3").Font.Bold = $true
    $ws.Range("K3").Font.Size = 14
    $ws.Range("K3").HorizontalAlignment = -4108
    $ws.Range("K3").Borders.Item(7).LineStyle = 1
    $ws.Range("K3").Borders.Item(8).LineStyle = 1
    $ws.Range("K3").Borders.Item(9).LineStyle = 1
    $ws.Range("K3").Borders.Item(10).LineStyle = 1

    # --- Data rows (row 4-13) : left block keyed off column A, right block
    #     keyed off column G (which retains its own, differently-sorted,
    #     flight order) ---
    for ($r = 4; $r -le 13; $r++) {
        $leftFlight = $ws.Cells.Item($r, 1).Value2
        $rightFlight = $ws.Cells.Item($r, 7).Value2

        $eCell = $ws.Range("E$r")
        $eCell.Value = $samples[$leftFlight]
        $eCell.NumberFormat = "General"

        $kCell = $ws.Range("K$r")
        $kCell.Value = $samples[$rightFlight]
        $kCell.NumberFormat = "General"
        $kCell.Font.Size = 14
        $kCell.HorizontalAlignment = -4108
        $kCell.Borders.Item(7).LineStyle = 1
        $kCell.Borders.Item(8).LineStyle = 1
        $kCell.Borders.Item(9).LineStyle = 1
        $kCell.Borders.Item(10).LineStyle = 1
    }

    # Re-establish the sort annotation over the shifted "sorted" block.
    $ws.Sort.SortFields.Clear()
    $ws.Sort.SortFields.Add($ws.Range("H4"))
    $ws.Sort.SetRange($ws.Range("G4:J13"))
    $ws.Sort.Header = 0
    $ws.Sort.Apply()

    # Size the two new columns to fit their content.
    $ws.Columns("E").AutoFit()
    $ws.Columns("K").AutoFit()
}

# --- Sheet-specific touch-ups -------------------------------------------------

$wsTropo = $wb.Worksheets.Item(1)
$wsTropo.Range("H17").Select()

$wsBl = $wb.Worksheets.Item(2)

# "bl tau" got an extra cosmetic pass: the box border around the "R^2"/
# "# UT Samples" pair (D:E) and "Mode"/"Mean" pair (G:H) was cleaned up so the
# two columns read as one continuous box instead of a double line.
$wsBl.Range("D3:D13").Borders.Item(10).LineStyle = 0
$wsBl.Range("H3:H13").Borders.Item(7).LineStyle = 0

$wsBl.Range("F3").Font.Bold = $true
$wsBl.Range("F3").Font.Size = 14
$wsBl.Range("F3").HorizontalAlignment = -4108

$wsBl.Range("F12").Select()
